# CEDS_Code_Change_Tracker.xlsx - add the new "Code Change Tracker" row
# documenting the commit:
#   "Added small changes to IO_functions.R and data_functions.R that were
#    not committed before the semester began."
#
# The table runs from row 3 through row 68 (row 1-2 are merged title /
# header rows). Row 67 is a normal, unshaded data row, so we clone its
# cell formatting (font/border/number-format/wrap, etc.) onto the new
# row 69 before filling in the new values - this keeps the new row's
# style identical to the rest of the table instead of inheriting the
# special shaded style used only by the last existing row (68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clone formatting for the new row from a representative existing row ---
$ws.Range("A67:H67").Copy()
$ws.Range("A69:H69").PasteSpecial(-4122)   # xlPasteFormats

# --- fill in the new row's content ---
$ws.Range("A69").Value = "Added small changes to IO_functions.R and data_functions.R that were not committed before the semester began."
$ws.Range("B69").Value = 65
$ws.Range("C69").Value = "Jon Seibert"
$ws.Range("D69").Value = "Committed"
$ws.Range("E69").Value = 42262
$ws.Range("F69").Value = "-"
$ws.Range("G69").Value = 42262
$ws.Range("H69").Value = "?"

# --- move the active selection to the cell after the new last row,
#     matching where Excel leaves the cursor after such an edit ---
$ws.Range("I69").Select()
